# Update "Revised Calcs_India" sheet: row 10/11 label + values change from
# "Peak Power Demand after Storage and DR[...] : MostRecentRun" to
# "Peak Hour Electricity Demand by Season[...] : NoSettings", with refreshed
# underlying numbers (new source run). Dependent formulas in DRC-BDRC and
# DRC-PADRC recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Revised Calcs_India")

$ws.Range("A10").Value = "Peak Hour Electricity Demand by Season[summer] : NoSettings"
$ws.Range("A11").Value = "Peak Hour Electricity Demand by Season[winter] : NoSettings"

$row10 = @(164581,167124,184567,199680,216498,232792,248850,264673,280541,302367,324258,346278,368229,390186,414236,438370,462533,486619,511017,535347,559541,584045,608681,633204,653169,675065,695867,715551,737951,757599,778490,799341)
$row11 = @(144301,144210,157047,167875,178025,187719,197150,206417,215696,226630,237552,248582,259703,270844,283077,295421,307771,320177,332695,346994,361314,375578,390049,404517,416752,428795,441039,453324,465684,477938,490292,502594)

for ($i = 0; $i -lt $row10.Length; $i++) {
    $ws.Cells.Item(10, 2 + $i).Value = $row10[$i]
    $ws.Cells.Item(11, 2 + $i).Value = $row11[$i]
}

# Selection changes recorded in the saved views
$ws.Range("C12").Select()

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B30").Select()

$wb.Save()
